$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EBAY")

# Row 14 - Gross Margin
$ws.Range("D14").Value = 0.757
$ws.Range("E14").Value = 0.7577
$ws.Range("F14").Value = 0.752
$ws.Range("G14").Value = 0.7527

# Row 15 - EBIT Margin
$ws.Range("B15").Value = 0.2633
$ws.Range("D15").Value = 0.2482
$ws.Range("E15").Value = 0.2323
$ws.Range("F15").Value = 0.2132
$ws.Range("G15").Value = 0.2155

# Row 16 - EBT margin
$ws.Range("D16").Value = 0.2859
$ws.Range("E16").Value = 0.246
$ws.Range("F16").Value = 0.2004
$ws.Range("G16").Value = 0.2025

# Row 17 - Net Profit Margin
$ws.Range("D17").Value = 0.5638
$ws.Range("E17").Value = 0.5569
$ws.Range("F17").Value = 0.5444
$ws.Range("G17").Value = 0.2068

# Row 18 - Free Cash Flow Margin
$ws.Range("B18").Value = 0.2105
$ws.Range("D18").Value = 0.2115
$ws.Range("E18").Value = 0.3292
$ws.Range("F18").Value = 0.316
$ws.Range("G18").Value = 0.3

# Row 29 - EBITDA Margin
$ws.Range("D29").Value = 0.3186
$ws.Range("E29").Value = 0.3068
$ws.Range("F29").Value = 0.2929
$ws.Range("G29").Value = 0.2883

# Row 30 - Operating Cash Flow Margin
$ws.Range("D30").Value = 0.2592
$ws.Range("E30").Value = 0.3742
$ws.Range("F30").Value = 0.3671
$ws.Range("G30").Value = 0.3606
